$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
  2  = @{ D = "47.346.22"; E = "  +2.67%  " }
  3  = @{ D = "2.504.22";  E = "  +2.10%  " }
  4  = @{ E = "  +0.14%  " }
  5  = @{ D = "323.98";    E = "  +0.44%  " }
  6  = @{ D = "109.44";    E = "  +3.80%  " }
  7  = @{ E = "  +1.40%  " }
  8  = @{ E = "  +0.05%  " }
  9  = @{ E = "  +0.70%  " }
  10 = @{ D = "38.96";     E = "  +7.61%  " }
  11 = @{ E = "  +1.16%  " }
  12 = @{ E = "  +1.03%  " }
  13 = @{ D = "18.54";     E = "  +0.34%  " }
  14 = @{ D = "7.22";      E = "  +1.87%  " }
  15 = @{ D = "2.894.30";  E = "  +2.43%  " }
  16 = @{ D = "2.503.35";  E = "  +2.12%  " }
  17 = @{ E = "  +1.77%  " }
  18 = @{ D = "47.303.43"; E = "  +3.00%  " }
  19 = @{ E = "  +3.42%  " }
  20 = @{ D = "6.73";      E = "  +4.47%  " }
  21 = @{ D = "0.0₃0947"; E = "  +1.43%  " }
  22 = @{ D = "71.07";     E = "  -0.95%  " }
  23 = @{ D = "2.58";      E = "  +8.59%  " }
  24 = @{ D = "249.81";    E = "  +0.69%  " }
  25 = @{ D = "2.60";      E = "  +3.42%  " }
  26 = @{ E = "  +0.43%  " }
  27 = @{ E = "  -0.02%  " }
  28 = @{ E = "  +4.90%  " }
  29 = @{ D = "10.04";     E = "  +3.60%  " }
  30 = @{ D = "35.89";     E = "  +6.58%  " }
  31 = @{ E = "  +5.00%  " }
  32 = @{ D = "50.06";     E = "  +1.22%  " }
  33 = @{ D = "19.99";     E = "  -2.35%  " }
  34 = @{ E = "  +3.25%  " }
  35 = @{ E = "  +4.08%  " }
  36 = @{ E = "  +0.23%  " }
  37 = @{ D = "2.00";      E = "  +5.03%  " }
  38 = @{ E = "  +3.84%  " }
  40 = @{ E = "  +1.34%  " }
  41 = @{ E = "  -1.76%  " }
  42 = @{ D = "121.99";    E = "  -4.08%  " }
  43 = @{ D = "21.41";     E = "  +2.46%  " }
  44 = @{ E = "  +2.14%  " }
  45 = @{ D = "1.990.87";  E = "  +1.38%  " }
  46 = @{ D = "3.05";      E = "  +2.28%  " }
  47 = @{ E = "  -1.37%  " }
  48 = @{ E = "  -3.05%  " }
  49 = @{ E = "  -1.29%  " }
  50 = @{ D = "5.32";      E = "  +8.67%  " }
  51 = @{ D = "78.46";     E = "  +0.81%  " }
}

# Rows whose new Price text is pure numeric-looking (e.g. "323.98") need the
# cell forced to Text format first, otherwise Excel auto-converts it to a
# number and loses the original "price displayed as text" formatting used
# throughout this sheet.
$numericLooking = @(5, 6, 10, 13, 14, 20, 22, 23, 24, 25, 29, 30, 32, 33, 37, 42, 43, 46, 50, 51)

foreach ($row in $updates.Keys) {
  $vals = $updates[$row]
  if ($vals.ContainsKey("D")) {
    $cell = $ws.Range("D$row")
    if ($numericLooking -contains $row) {
      $cell.NumberFormat = "@"
    }
    $cell.Value = $vals["D"]
  }
  if ($vals.ContainsKey("E")) {
    $ws.Range("E$row").Value = $vals["E"]
  }
}
